# The edit removes the subscript-base letter "b" that directly follows each
# "log" in the two logarithm-identity cells of the LLR table
# ("Formula of Product" row -> LLR_L005, "Formula of Ratio" row -> LLR_L006),
# turning "logb(...)" into "log(...)" throughout both formulas:
#
#   " logb(xy)=logb(x)+logb(y)."      ->  " log(xy)=log(x)+log(y)."
#   " logb(x/y) =logb(x)-logb(y)."    ->  " log(x/y) =log(x)-log(y)."
#
# In the underlying XML, each "log" is its own bold run and the following
# "b(...)" text is a separate, non-bold run. Deleting just the leading "b"
# character of each such run (rather than doing a whole-cell text replace)
# preserves every run's existing character formatting exactly, instead of
# collapsing the cell into one differently-formatted run.

$d = $word.ActiveDocument

function Find-DescriptionCellByLabel($labelText) {
    # The requirement tables have 3 columns: Code | Requirements | Description.
    # Locate the row whose "Requirements" (column 2) cell starts with
    # $labelText and return its "Description" (column 3) cell.
    foreach ($tbl in $d.Tables) {
        for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
            $labelCell = $tbl.Cell($r, 2)
            if ($labelCell.Range.Text.StartsWith($labelText)) {
                return $tbl.Cell($r, 3)
            }
        }
    }
    return $null
}

function Remove-LogSubscriptB($cell) {
    $cellStart = $cell.Range.Start
    $cellEnd = $cell.Range.End

    # Walk forward through the cell looking for "log"; immediately after
    # each occurrence, delete the single "b" character that starts the
    # following (non-bold) run, if present.
    $pos = $cellStart
    while ($true) {
        $fr = $d.Range($pos, $cellEnd)
        $found = $fr.Find.Execute("log")
        if (-not $found) { break }

        $afterLog = $fr.End
        $bChar = $d.Range($afterLog, $afterLog + 1)
        if ($bChar.Text -eq "b") {
            $bChar.Delete()
            # Cell end shrinks by one character since we deleted from inside it.
            $cellEnd = $cellEnd - 1
        }
        $pos = $afterLog
    }
}

Remove-LogSubscriptB(Find-DescriptionCellByLabel("Formula of Product"))
Remove-LogSubscriptB(Find-DescriptionCellByLabel("Formula of Ratio"))

# Re-look-up the cells fresh (post-edit) for the status output below, since
# a cell reference captured before an edit can report stale cached text.
Write-Output ("Formula of Product -> [" + (Find-DescriptionCellByLabel("Formula of Product")).Range.Text + "]")
Write-Output ("Formula of Ratio -> [" + (Find-DescriptionCellByLabel("Formula of Ratio")).Range.Text + "]")
